$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2025-05-17T18:39:20.262741"
$ws.Range("B3").Value = "MXN"
$ws.Range("C3").Value = "USD"
$ws.Range("D3").Value = 500
$ws.Range("E3").Value = 0.05137109451253968
$ws.Range("F3").Value = 25.68554725626984
